$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (B, C, D, E) keep their values as literal text,
# matching the inlineStr cell type in the original workbook (avoids Excel
# auto-converting values like "61.690.96" or "1.00" into numbers/dates).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.690.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.457.01'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.17%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.71'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.04'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +7.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.456.63'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.475'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.70'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.126'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.392'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.047.99'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.88'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +9.55%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.482.18'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.00%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000175'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.825.38'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.25'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +8.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.39'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.57'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.43'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.565'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.39'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.43%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000122'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.606.28'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.58%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.76'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.89%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.50'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -9.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.22'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.19'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.21%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '24.23'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.489.66'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.01'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.34%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.58'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.16'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '166.61'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0786'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '27.05'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +6.34%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.74%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.51'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.86%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '42.55'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.73'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.17'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.564.69'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.93'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.06%  '
